# Add the new "SOT-E2" worksheet (E2 SKU pricing) as a copy of "SOT-N2",
# positioned after "global license pricing" (i.e. as the last tab), then
# replace its contents with the E2 pricing figures.

$wb = $excel.ActiveWorkbook

$src  = $wb.Worksheets.Item("SOT-N2")
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $last)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "SOT-E2"

# The E2 sheet doesn't use the "Ext_RAM usd" (column B) figures per-row, so
# clear the data cells while keeping the column B header in row 1.
$ws.Range("B2:B24").ClearContents()

# New column headers for the on-demand / 1yr / 3yr custom-core E2 pricing.
$ws.Range("F1").Value = "E2_custom_core_on-deman"
$ws.Range("G1").Value = "E2_custom_core_1yr"
$ws.Range("H1").Value = "E2_custom_core_3yr"

# Region (col A) stays as copied from SOT-N2; fill in the E2 pricing data
# for columns C (custom_RAM no-commit) .. H (E2_custom_core_3yr).
$data = @(
    @(2, 0.00338514, 0.0021326380000000001, 0.0015233130000000001, 0.025255380000000001, 0.015910889000000001, 0.011364921),
    @(3, 0.0029235300000000001, 0.001841823, 0.0013155879999999999, 0.021811589999999999, 0.013741300999999999, 0.0098152150000000004),
    @(4, 0.00321609, 0.002026136, 0.00144724, 0.02399337, 0.015115823, 0.010797016),
    @(5, 0.0032181599999999999, 0.0020274400000000001, 0.001448172, 0.024016139999999998, 0.015130167999999999, 0.010807262999999999),
    @(6, 0.0037660200000000001, 0.0023725920000000002, 0.001694709, 0.028103699999999999, 0.017705331000000001, 0.012646665),
    @(7, 0.0040903199999999997, 0.0025769009999999999, 0.001840644, 0.03051939, 0.019227214999999999, 0.013733725),
    @(8, 0.0039304739999999998, 0.0024761980000000002, 0.0017687130000000001, 0.029331146999999998, 0.018478628, 0.013199017),
    @(9, 0.0037391099999999999, 0.0023556390000000001, 0.0016825989999999999, 0.02802642, 0.017656643999999999, 0.012611888999999999),
    @(10, 0.0032919899999999998, 0.0020739529999999999, 0.0014813949999999999, 0.024567450000000001, 0.015477493, 0.011055351999999999),
    @(11, 0.0037660200000000001, 0.0023725920000000002, 0.001694709, 0.028103699999999999, 0.017705331000000001, 0.012646665),
    @(12, 0.00351141, 0.0022121879999999999, 0.001580134, 0.026199299999999998, 0.016505559, 0.011789684999999999),
    @(13, 0.0032181599999999999, 0.0020274400000000001, 0.001448172, 0.024013380000000001, 0.015128429000000001, 0.010806021000000001),
    @(14, 0.0035107200000000002, 0.0022117529999999999, 0.0015798240000000001, 0.026199299999999998, 0.016505559, 0.011789684999999999),
    @(15, 0.0032181599999999999, 0.0020274400000000001, 0.001448172, 0.024013380000000001, 0.015128429000000001, 0.010806021000000001),
    @(16, 0.0037391099999999999, 0.0023556390000000001, 0.0016825989999999999, 0.02802642, 0.017656643999999999, 0.012611888999999999),
    @(17, 0.00351141, 0.0022121879999999999, 0.001580134, 0.026199299999999998, 0.016505559, 0.011789684999999999),
    @(18, 0.0046402500000000003, 0.0029233570000000001, 0.002088112, 0.034624200000000001, 0.021813246000000001, 0.01558089),
    @(19, 0.0037391099999999999, 0.0023556390000000001, 0.0016825989999999999, 0.02802642, 0.017656643999999999, 0.012611888999999999),
    @(20, 0.0036059400000000002, 0.0022717420000000002, 0.001622673, 0.026909309999999999, 0.016952865000000001, 0.012109188999999999),
    @(21, 0.0041475899999999996, 0.0026129809999999999, 0.001866415, 0.030950640000000001, 0.019498903000000001, 0.013927788),
    @(22, 0.0032919899999999998, 0.0020739529999999999, 0.0014813949999999999, 0.024567450000000001, 0.015477493, 0.011055351999999999),
    @(23, 0.0037660200000000001, 0.0023725920000000002, 0.001694709, 0.028103699999999999, 0.017705331000000001, 0.012646665),
    @(24, 0.0040903199999999997, 0.0025769009999999999, 0.001840644, 0.03051939, 0.019227214999999999, 0.013733725)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
}

# Match the saved selection/active-cell state on the new sheet and make it
# the active tab (mirrors the source file no longer being tabSelected).
$ws.Activate()
$ws.Range("G12").Select()
